$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column F ("Categoría"), shifting
# Categoría..Observaciones (old F:Y) to the right (new G:Z).
$ws.Columns("F:F").Insert()

# Populate the new "Sede de Registro" column.
$ws.Range("F1").Value = "Sede de Registro"
$ws.Range("F2").Value = "PUNO"
$ws.Range("F3").Value = "AREQUIPA"
$ws.Range("F4").Value = "CUSCO"
$ws.Range("F5").Value = "INVALID_SEDE"
$ws.Range("F6").Value = "LIMA"
